# ECS reading refactor, moved to separate class
#
# Rename the "envs" sheet to "Parameters" and replace its 3-column
# param/env01/env02 layout with a simple 2-column param/value layout,
# keeping only a single "ECS sheet" / "5.1ECS" row of data and clearing
# the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet "envs" -> "Parameters"
$ws.Name = "Parameters"

# Delete the whole column C (env02), shifting nothing else - column B
# (env01) becomes the sole "value" column.
$ws.Columns.Item(3).Delete()

# Update header row
$ws.Range("A1").Value = "param"
$ws.Range("B1").Value = "value"

# Update data row 2, and clear out rows 3-5 (previously sheet/range_start/range_end)
$ws.Range("B2").Value = "5.1ECS"
$ws.Range("A2").Value = "ECS sheet"

$ws.Range("A3:B5").ClearContents()

# Adjust column A width to match the new narrower layout (ColumnWidth is in
# character units; the stored XML width includes the standard +5/7 padding,
# so back that off here to land on an exact "13" in the saved file).
$ws.Columns.Item(1).ColumnWidth = 12.2857142857143

# Move the selection to A3, matching the saved view state
$ws.Range("A3").Select()
